$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.73885133333334
$ws.Range("H2").Value = 107.216554
$ws.Range("I2").Value = 0.01949729408921566
$ws.Range("J2").Value = 0.01949729408921566
$ws.Range("M2").Value = 0.1419263333333333
$ws.Range("N2").Value = 0.425779
$ws.Range("O2").Value = 0.002583058778296354
$ws.Range("P2").Value = 0.002583058778296354
$ws.Range("Q2").Value = 5.072284127285112
$ws.Range("R2").Value = 45.650557145566
$ws.Range("S2").Value = 0.00005036265665017413
$ws.Range("T2").Value = 0.00005036265665017411
$ws.Range("G3").Value = 35.73885133333334
$ws.Range("H3").Value = 107.216554
$ws.Range("I3").Value = 0.01949729408921566
$ws.Range("J3").Value = 0.01949729408921566
$ws.Range("O3").Value = 0.001399682868699959
$ws.Range("P3").Value = 0.001399682868699959
$ws.Range("Q3").Value = 2.748520187690889
$ws.Range("R3").Value = 24.736681689218
$ws.Range("S3").Value = 0.00002729002852268013
$ws.Range("T3").Value = 0.00002729002852268013
$ws.Range("G4").Value = 35.73885133333334
$ws.Range("H4").Value = 107.216554
$ws.Range("I4").Value = 0.01949729408921566
$ws.Range("J4").Value = 0.01949729408921566
$ws.Range("M4").Value = 2.613991
$ws.Range("N4").Value = 7.841973
$ws.Range("O4").Value = 0.04757462720522382
$ws.Range("P4").Value = 0.04757462720522382
$ws.Range("Q4").Value = 93.42103573567134
$ws.Range("R4").Value = 840.7893216210421
$ws.Range("S4").Value = 0.000927576497805049
$ws.Range("T4").Value = 0.0009275764978050488
$ws.Range("G5").Value = 35.73885133333334
$ws.Range("H5").Value = 107.216554
$ws.Range("I5").Value = 0.01949729408921566
$ws.Range("J5").Value = 0.01949729408921566
$ws.Range("M5").Value = 52.11224233333333
$ws.Range("N5").Value = 156.336727
$ws.Range("O5").Value = 0.9484426311477799
$ws.Range("P5").Value = 0.9484426311477798
$ws.Range("Q5").Value = 1862.43168139764
$ws.Range("R5").Value = 16761.88513257876
$ws.Range("S5").Value = 0.01849206490623776
$ws.Range("T5").Value = 0.01849206490623775
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.9215900675332435
$ws.Range("J6").Value = 0.9215900675332435
$ws.Range("M6").Value = 0.1419263333333333
$ws.Range("N6").Value = 0.425779
$ws.Range("O6").Value = 0.002583058778296354
$ws.Range("P6").Value = 0.002583058778296354
$ws.Range("Q6").Value = 239.7546372344089
$ws.Range("R6").Value = 2157.79173510968
$ws.Range("S6").Value = 0.002380521313932474
$ws.Range("T6").Value = 0.002380521313932474
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.9215900675332435
$ws.Range("J7").Value = 0.9215900675332435
$ws.Range("O7").Value = 0.001399682868699959
$ws.Range("P7").Value = 0.001399682868699959
$ws.Range("S7").Value = 0.001289933829490319
$ws.Range("T7").Value = 0.001289933829490319
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.9215900675332435
$ws.Range("J8").Value = 0.9215900675332435
$ws.Range("M8").Value = 2.613991
$ws.Range("N8").Value = 7.841973
$ws.Range("O8").Value = 0.04757462720522382
$ws.Range("P8").Value = 0.04757462720522382
$ws.Range("Q8").Value = 4415.787044022906
$ws.Range("R8").Value = 39742.08339620616
$ws.Range("S8").Value = 0.0438443038989311
$ws.Range("T8").Value = 0.0438443038989311
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.9215900675332435
$ws.Range("J9").Value = 0.9215900675332435
$ws.Range("M9").Value = 52.11224233333333
$ws.Range("N9").Value = 156.336727
$ws.Range("O9").Value = 0.9484426311477799
$ws.Range("P9").Value = 0.9484426311477798
$ws.Range("Q9").Value = 88032.65372012198
$ws.Range("R9").Value = 792293.8834810978
$ws.Range("S9").Value = 0.8740753084908897
$ws.Range("T9").Value = 0.8740753084908895
$ws.Range("G10").Value = 93.641553
$ws.Range("H10").Value = 280.924659
$ws.Range("I10").Value = 0.05108605424341119
$ws.Range("J10").Value = 0.05108605424341119
$ws.Range("M10").Value = 0.1419263333333333
$ws.Range("N10").Value = 0.425779
$ws.Range("O10").Value = 0.002583058778296354
$ws.Range("P10").Value = 0.002583058778296354
$ws.Range("Q10").Value = 13.290202264929
$ws.Range("R10").Value = 119.611820384361
$ws.Range("S10").Value = 0.000131958280861967
$ws.Range("T10").Value = 0.000131958280861967
$ws.Range("G11").Value = 93.641553
$ws.Range("H11").Value = 280.924659
$ws.Range("I11").Value = 0.05108605424341119
$ws.Range("J11").Value = 0.05108605424341119
$ws.Range("O11").Value = 0.001399682868699959
$ws.Range("P11").Value = 0.001399682868699959
$ws.Range("Q11").Value = 7.201566061167
$ws.Range("R11").Value = 64.814094550503
$ws.Range("S11").Value = 0.00007150427495397948
$ws.Range("T11").Value = 0.00007150427495397948
$ws.Range("G12").Value = 93.641553
$ws.Range("H12").Value = 280.924659
$ws.Range("I12").Value = 0.05108605424341119
$ws.Range("J12").Value = 0.05108605424341119
$ws.Range("M12").Value = 2.613991
$ws.Range("N12").Value = 7.841973
$ws.Range("O12").Value = 0.04757462720522382
$ws.Range("P12").Value = 0.04757462720522382
$ws.Range("Q12").Value = 244.778176768023
$ws.Range("R12").Value = 2203.003590912207
$ws.Range("S12").Value = 0.00243039998601613
$ws.Range("T12").Value = 0.00243039998601613
$ws.Range("G13").Value = 93.641553
$ws.Range("H13").Value = 280.924659
$ws.Range("I13").Value = 0.05108605424341119
$ws.Range("J13").Value = 0.05108605424341119
$ws.Range("M13").Value = 52.11224233333333
$ws.Range("N13").Value = 156.336727
$ws.Range("O13").Value = 0.9484426311477799
$ws.Range("P13").Value = 0.9484426311477798
$ws.Range("Q13").Value = 4879.871302405677
$ws.Range("R13").Value = 43918.84172165109
$ws.Range("S13").Value = 0.04845219170157911
$ws.Range("T13").Value = 0.04845219170157911
$ws.Range("G14").Value = 14.34625366666667
$ws.Range("H14").Value = 43.038761
$ws.Range("I14").Value = 0.007826584134129748
$ws.Range("J14").Value = 0.007826584134129748
$ws.Range("M14").Value = 0.1419263333333333
$ws.Range("N14").Value = 0.425779
$ws.Range("O14").Value = 0.002583058778296354
$ws.Range("P14").Value = 0.002583058778296354
$ws.Range("Q14").Value = 2.036111179979889
$ws.Range("R14").Value = 18.325000619819
$ws.Range("S14").Value = 0.00002021652685173881
$ws.Range("T14").Value = 0.00002021652685173881
$ws.Range("G15").Value = 14.34625366666667
$ws.Range("H15").Value = 43.038761
$ws.Range("I15").Value = 0.007826584134129748
$ws.Range("J15").Value = 0.007826584134129748
$ws.Range("O15").Value = 0.001399682868699959
$ws.Range("P15").Value = 0.001399682868699959
$ws.Range("Q15").Value = 1.103308202404111
$ws.Range("R15").Value = 9.929773821637001
$ws.Range("S15").Value = 0.00001095473573298031
$ws.Range("T15").Value = 0.00001095473573298031
$ws.Range("G16").Value = 14.34625366666667
$ws.Range("H16").Value = 43.038761
$ws.Range("I16").Value = 0.007826584134129748
$ws.Range("J16").Value = 0.007826584134129748
$ws.Range("M16").Value = 2.613991
$ws.Range("N16").Value = 7.841973
$ws.Range("O16").Value = 0.04757462720522382
$ws.Range("P16").Value = 0.04757462720522382
$ws.Range("Q16").Value = 37.50097796838367
$ws.Range("R16").Value = 337.508801715453
$ws.Range("S16").Value = 0.0003723468224715422
$ws.Range("T16").Value = 0.0003723468224715422
$ws.Range("G17").Value = 14.34625366666667
$ws.Range("H17").Value = 43.038761
$ws.Range("I17").Value = 0.007826584134129748
$ws.Range("J17").Value = 0.007826584134129748
$ws.Range("M17").Value = 52.11224233333333
$ws.Range("N17").Value = 156.336727
$ws.Range("O17").Value = 0.9484426311477799
$ws.Range("P17").Value = 0.9484426311477798
$ws.Range("Q17").Value = 747.6154476528053
$ws.Range("R17").Value = 6728.539028875247
$ws.Range("S17").Value = 0.007423066049073487
$ws.Range("T17").Value = 0.007423066049073486
